$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.250.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.95%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.314.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.53%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.79%  "

# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.311.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.49%  "

# Row 9
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.22%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.835.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.05%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.117"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.96%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.305.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.64%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.12%  "

# Row 20
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.063.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.965"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.88%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.25%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.75%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.29%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.18%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "584.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.73%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.148"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.84%  "

# Row 40
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.53%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.17%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.106.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.84%  "

# Row 45
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.96%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0400"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.127"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.90%  "
